# The source data files for this Transition-Matrix workbook were moved
# into a "Transition Matrices" subfolder, and the "ScreenRecStarted"
# gaze-state label was renamed to "0_unstated" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename every "ScreenRecStarted" occurrence (header + transition-pair
# labels in column A) to "0_unstated".
$ws.Range("G1").Value = "0_unstated"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Update the active selection to G1 (matches the saved view state).
$ws.Range("G1").Select()
